$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo and update feature names / score type labels (rows 4-8)
$ws.Range("A4").Value = "HAS-BLED Score"
$ws.Range("E4").Value = "cont"

$ws.Range("A5").Value = "High-Risk Alcohol Consumption"
$ws.Range("E5").Value = "cat"

$ws.Range("A6").Value = "Platelet Aggregation Inhibitor Therapy"
$ws.Range("E6").Value = "cat"

$ws.Range("A7").Value = "Oral Anticoagulation Therapy"
$ws.Range("E7").Value = "cat"

$ws.Range("A8").Value = "Perioperative Bridging Therapy"
$ws.Range("E8").Value = "cat"

# Update the selected range shown in the sheet view
$ws.Range("A9:G25").Select()

$wb.Save()
